# Generate Report for Handoff
# Updates the "b.md" row across the Overview, zh-cn and de-de sheets to
# reflect that a new handoff has been generated (status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# zh-cn/de-de sheets record the new handoff file + timestamp).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-11 02:47:08"

# --- de-de sheet ------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-11 02:47:15"
